$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (the "Förändrad" date column) for rows 2 through 18 is being
# updated from 45184 (2023-09-15) to 45185 (2023-09-16).
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
